# This script updates the "dSF" (column F) values on Sheet1 to match a
# repulled/recomputed dataset, per commit message: "repull data, push all
# data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> new value for column F (dSF)
$updates = @{
    2  = 0
    6  = 2
    7  = 3
    10 = 1
    11 = 3
    15 = 0
    24 = -3
    25 = -6
    26 = -12
    30 = -5
    34 = -2
    38 = -4
    39 = 10
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
